$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
